$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 36; existing rows 36:155 shift down to 37:156.
$ws.Rows("36:36").Insert()

# Populate the newly inserted row 36 with the new weekly data point.
$ws.Cells.Item(36, 1).Value = 3
$ws.Cells.Item(36, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(36, 3).Value = "Coquimbo"
$ws.Cells.Item(36, 4).Value = 44481
$ws.Cells.Item(36, 5).Value = 5
$ws.Cells.Item(36, 6).Value = 100112001
$ws.Cells.Item(36, 7).Value = "Berenjena"
$ws.Cells.Item(36, 8).Value = "Sin especificar"
$ws.Cells.Item(36, 9).Value = "Primera"
$ws.Cells.Item(36, 10).Value = 30
$ws.Cells.Item(36, 11).Value = 9000
$ws.Cells.Item(36, 12).Value = 9000
$ws.Cells.Item(36, 13).Value = 9000
$ws.Cells.Item(36, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(36, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(36, 16).Value = 150
$ws.Cells.Item(36, 17).Value = 60
$ws.Cells.Item(36, 18).Value = "Hortaliza"
